$wb = $excel.ActiveWorkbook

# Rename two worksheet tabs ("coeff" -> "coef")
$wsAng1 = $wb.Worksheets.Item("sensitivity coeff., ang1")
$wsAng1.Name = "sensitivity coef, ang1"

$wsAng2 = $wb.Worksheets.Item("sensitivity coeff, ang2")
$wsAng2.Name = "sensitivity coef, ang2"

# Update cell selections on a couple of sheets to reflect where the user last clicked
$wsKpp = $wb.Worksheets.Item("KPP, no WFE")
$wsKpp.Activate() | Out-Null
$wsKpp.Range("C3").Select() | Out-Null

$wsAng1.Activate() | Out-Null
$wsAng1.Range("H40").Select() | Out-Null

# Make the Glossary sheet the active (selected) tab
$wsGlossary = $wb.Worksheets.Item("Glossary")
$wsGlossary.Activate() | Out-Null
